$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Apply the bold/bordered style (from A3) to the new A-column index cells A10:A19
$ws.Cells.Item(3, 1).Copy()
$ws.Range("A10:A19").PasteSpecial(-4122)

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "Gaussian-Quadrature"
$ws.Cells.Item(10, 3).Value = 1.193525908319875
$ws.Cells.Item(10, 4).Value = 0.8636495835432553
$ws.Cells.Item(10, 5).Value = 0.9527511789956422
$ws.Cells.Item(10, 6).Value = 1.193525908319875
$ws.Cells.Item(10, 7).Value = 0.8947636091639845
$ws.Cells.Item(10, 8).Value = 1.160794165004746
$ws.Cells.Item(10, 9).Value = 1.008482133597503
$ws.Cells.Item(10, 10).Value = 0.8636495835432553
$ws.Cells.Item(10, 11).Value = 0.9082003812694488
$ws.Cells.Item(10, 12).Value = 1.050863144794662
$ws.Cells.Item(10, 13).Value = 1.012327763104168

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "Spiral-90deg-10rot-5space"
$ws.Cells.Item(11, 3).Value = 0.758057192404886
$ws.Cells.Item(11, 4).Value = 1.230266664069986
$ws.Cells.Item(11, 5).Value = 1.104538688544797
$ws.Cells.Item(11, 6).Value = 0.758057192404886
$ws.Cells.Item(11, 7).Value = 0.8182279735178831
$ws.Cells.Item(11, 8).Value = 1.806161283372751
$ws.Cells.Item(11, 9).Value = 0.9865825605729284
$ws.Cells.Item(11, 10).Value = 1.230266664069986
$ws.Cells.Item(11, 11).Value = 1.167402676307392
$ws.Cells.Item(11, 12).Value = 0.9627299343561389
$ws.Cells.Item(11, 13).Value = 1.117305727080539

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "Spiral-90deg-15rot-5space"
$ws.Cells.Item(12, 3).Value = 0.7589811257689935
$ws.Cells.Item(12, 4).Value = 1.234181831841632
$ws.Cells.Item(12, 5).Value = 1.102941784976655
$ws.Cells.Item(12, 6).Value = 0.7589811257689935
$ws.Cells.Item(12, 7).Value = 0.8207279325533646
$ws.Cells.Item(12, 8).Value = 1.799472068638838
$ws.Cells.Item(12, 9).Value = 0.9855191801030324
$ws.Cells.Item(12, 10).Value = 1.234181831841632
$ws.Cells.Item(12, 11).Value = 1.168561808409143
$ws.Cells.Item(12, 12).Value = 0.9637714670890684
$ws.Cells.Item(12, 13).Value = 1.116970653980419

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "Spiral-90deg-10rot-3space"
$ws.Cells.Item(13, 3).Value = 0.7577226794572892
$ws.Cells.Item(13, 4).Value = 1.231714785333718
$ws.Cells.Item(13, 5).Value = 1.10393165149707
$ws.Cells.Item(13, 6).Value = 0.7577226794572892
$ws.Cells.Item(13, 7).Value = 0.8190883809577582
$ws.Cells.Item(13, 8).Value = 1.806415790952808
$ws.Cells.Item(13, 9).Value = 0.9857676872138565
$ws.Cells.Item(13, 10).Value = 1.231714785333718
$ws.Cells.Item(13, 11).Value = 1.167823218415394
$ws.Cells.Item(13, 12).Value = 0.9627729489363418
$ws.Cells.Item(13, 13).Value = 1.11744016256875

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "NoRotation-tilt60deg"
$ws.Cells.Item(14, 3).Value = 0.5632000000000011
$ws.Cells.Item(14, 4).Value = 0.9345080000000023
$ws.Cells.Item(14, 5).Value = 1.373200000000002
$ws.Cells.Item(14, 6).Value = 0.5632000000000011
$ws.Cells.Item(14, 7).Value = 0.7573960000000016
$ws.Cells.Item(14, 8).Value = 1.487591999999998
$ws.Cells.Item(14, 9).Value = 1.151316000000002
$ws.Cells.Item(14, 10).Value = 0.9345080000000023
$ws.Cells.Item(14, 11).Value = 1.153854000000002
$ws.Cells.Item(14, 12).Value = 0.8585270000000015
$ws.Cells.Item(14, 13).Value = 1.044535333333334

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "Rotation-NoTilt"
$ws.Cells.Item(15, 3).Value = 0.42
$ws.Cells.Item(15, 4).Value = 0.14
$ws.Cells.Item(15, 5).Value = 1.7604375
$ws.Cells.Item(15, 6).Value = 0.42
$ws.Cells.Item(15, 7).Value = 0.45
$ws.Cells.Item(15, 8).Value = 1.6001375
$ws.Cells.Item(15, 9).Value = 1.411399999999998
$ws.Cells.Item(15, 10).Value = 0.14
$ws.Cells.Item(15, 11).Value = 0.9502187500000001
$ws.Cells.Item(15, 12).Value = 0.685109375
$ws.Cells.Item(15, 13).Value = 0.9636624999999999

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "Rotation-60detTilt"
$ws.Cells.Item(16, 3).Value = 0.6664054890496031
$ws.Cells.Item(16, 4).Value = 0.4883927640064026
$ws.Cells.Item(16, 5).Value = 1.426963998924803
$ws.Cells.Item(16, 6).Value = 0.6664054890496031
$ws.Cells.Item(16, 7).Value = 0.6838353611775988
$ws.Cells.Item(16, 8).Value = 1.350804546355197
$ws.Cells.Item(16, 9).Value = 1.225006526976004
$ws.Cells.Item(16, 10).Value = 0.4883927640064026
$ws.Cells.Item(16, 11).Value = 0.9576783814656027
$ws.Cells.Item(16, 12).Value = 0.8120419352576029
$ws.Cells.Item(16, 13).Value = 0.9735681144149346

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "HexGrid-90degTilt5degRes"
$ws.Cells.Item(17, 3).Value = 0.9877470517650069
$ws.Cells.Item(17, 4).Value = 0.9957049358836119
$ws.Cells.Item(17, 5).Value = 0.995483022488447
$ws.Cells.Item(17, 6).Value = 0.9877470517650069
$ws.Cells.Item(17, 7).Value = 0.9912090312430195
$ws.Cells.Item(17, 8).Value = 0.9983375161198511
$ws.Cells.Item(17, 9).Value = 0.9916453569570173
$ws.Cells.Item(17, 10).Value = 0.9957049358836119
$ws.Cells.Item(17, 11).Value = 0.9955939791860294
$ws.Cells.Item(17, 12).Value = 0.9916705154755181
$ws.Cells.Item(17, 13).Value = 0.9933544857428256

$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = "HexGrid-90degTilt22p5degRes"
$ws.Cells.Item(18, 3).Value = 1.054455857472574
$ws.Cells.Item(18, 4).Value = 1.030516624576193
$ws.Cells.Item(18, 5).Value = 0.9632128425060728
$ws.Cells.Item(18, 6).Value = 1.054455857472574
$ws.Cells.Item(18, 7).Value = 0.995467729785534
$ws.Cells.Item(18, 8).Value = 0.9308028245614358
$ws.Cells.Item(18, 9).Value = 0.9854158957792784
$ws.Cells.Item(18, 10).Value = 1.030516624576193
$ws.Cells.Item(18, 11).Value = 0.9968647335411329
$ws.Cells.Item(18, 12).Value = 1.025660295506853
$ws.Cells.Item(18, 13).Value = 0.9933119624468479

$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(19, 3).Value = 0.9615077257422631
$ws.Cells.Item(19, 4).Value = 1.029101065489702
$ws.Cells.Item(19, 5).Value = 0.9790213677677105
$ws.Cells.Item(19, 6).Value = 0.9615077257422631
$ws.Cells.Item(19, 7).Value = 1.049377421082012
$ws.Cells.Item(19, 8).Value = 0.9107184732084199
$ws.Cells.Item(19, 9).Value = 0.9698016854779219
$ws.Cells.Item(19, 10).Value = 1.029101065489702
$ws.Cells.Item(19, 11).Value = 1.004061216628706
$ws.Cells.Item(19, 12).Value = 0.9827844711854847
$ws.Cells.Item(19, 13).Value = 0.9832546231280049

$excel.CutCopyMode = $false
